$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Insert a new worksheet "Sheet2" right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# 2. Populate Sheet2's header row by copying formatted cells from Sheet1
#    (keeps the exact same direct-format styles Sheet1 uses).
$ws1.Range("A1").Copy($ws2.Range("A1"))   # "Title"
$ws1.Range("B1").Copy($ws2.Range("B1"))   # "Id"

# C2 is still styled "s=4" at this point (it only gets restyled to s=1 at the
# end of this script) -- borrow it now, while it is still the only cell on
# Sheet1 carrying that style, for every Sheet2 cell that needs it.
$ws1.Range("C2").Copy($ws2.Range("C1"))
$ws2.Range("C1").Value = "color"

$ws1.Range("C2").Copy($ws2.Range("A2"))   # "circle"
$ws1.Range("C2").Copy($ws2.Range("A3"))
$ws1.Range("C2").Copy($ws2.Range("A4"))
$ws1.Range("C2").Copy($ws2.Range("A5"))
$ws1.Range("C2").Copy($ws2.Range("A6"))
$ws1.Range("C2").Copy($ws2.Range("A7"))

# Numeric Id column: copy styled cells straight across (values already 1..6).
$ws1.Range("B2").Copy($ws2.Range("B2"))
$ws1.Range("B3").Copy($ws2.Range("B3"))
$ws1.Range("B4").Copy($ws2.Range("B4"))
$ws1.Range("B5").Copy($ws2.Range("B5"))
$ws1.Range("B6").Copy($ws2.Range("B6"))
$ws1.Range("B7").Copy($ws2.Range("B7"))

# Color column: same cycle of values Sheet1 keeps in column A (Red/Green/Blue).
$ws2.Range("C2").Value = "Red"
$ws2.Range("C3").Value = "Green"
$ws2.Range("C4").Value = "Blue"
$ws2.Range("C5").Value = "Red"
$ws2.Range("C6").Value = "Green"
$ws2.Range("C7").Value = "Blue"

# 3. Sheet1: C2 picks up the style (and, incidentally, the identical "circle"
#    value) that C3 already carries, landing on style s=1 exactly like the
#    rest of column C below the header.
$ws1.Range("C3").Copy($ws1.Range("C2"))
